# RSD Retrofit: Use LimType FX with NCAP_AF
#
# Inserts a new "LimType" column (with value "FX" on every data row) into
# the rsd_rtft table on the RSD_RTFT worksheet, between the existing
# "Attribute" and "Year" columns. Everything from the old "Year" column
# onward shifts one column to the right (Year, Pset_PN, IE, National).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RSD_RTFT")
$lo = $ws.ListObjects.Item(1)

# Capture the table's current data (header row + body rows) before touching
# anything, so we can rebuild it in the new column layout.
$headerRange = $lo.HeaderRowRange
$bodyRange = $lo.DataBodyRange
$numDataRows = $bodyRange.Rows.Count

$oldRange = $ws.Range($headerRange.Cells.Item(1, 1), $bodyRange.Cells.Item($numDataRows, $bodyRange.Columns.Count))
$oldData = $oldRange.Value()

$rows = $numDataRows + 1
$cols = 7
$newData = New-Object 'object[,]' $rows, $cols

# New header row: TimeSlice, Attribute, LimType, Year, Pset_PN, IE, National
$newData[0, 0] = "TimeSlice"
$newData[0, 1] = "Attribute"
$newData[0, 2] = "LimType"
$newData[0, 3] = "Year"
$newData[0, 4] = "Pset_PN"
$newData[0, 5] = "IE"
$newData[0, 6] = "National"

for ($r = 1; $r -lt $rows; $r++) {
    $srcR = $r + 1
    $newData[$r, 0] = $oldData[$srcR, 1]
    $newData[$r, 1] = $oldData[$srcR, 2]
    $newData[$r, 2] = "FX"
    $newData[$r, 3] = $oldData[$srcR, 3]
    $newData[$r, 4] = $oldData[$srcR, 4]
    $newData[$r, 5] = $oldData[$srcR, 5]
    $newData[$r, 6] = $oldData[$srcR, 6]
}

# Grow the table by one column (to the right) then write the rebuilt data,
# which also renames/reorders the ListColumns via their header cells.
$newTableRange = $ws.Range($headerRange.Cells.Item(1, 1), $bodyRange.Cells.Item($numDataRows, $bodyRange.Columns.Count + 1))
$lo.Resize($newTableRange)
$newTableRange.Value = $newData

# Best-fit the new Year column's width (now column E).
$ws.Columns("E:E").AutoFit()

# Make RSD_RTFT the active sheet/tab and leave the selection where the
# author left it.
$ws.Activate()
$ws.Range("J6").Select()
